# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The sheet currently spans A1:AC50 (player roster / stats). This adds
# three new trailing columns, AD (Wins), AE (Losses) and AF (Ties), with
# a header row styled like the other headers and the season record
# (64 wins, 98 losses, 0 ties) repeated down every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row --------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header cells (bold, centered,
# bordered) by copying the style from an existing header cell (A1).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
$wins = 64
$losses = 98
$ties = 0

for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
